# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets, matching the refreshed data pulled at commit 7921097.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F4"  = 11030
    "F5"  = 10205
    "F6"  = 591
    "F13" = 9580
    "F15" = 2435
    "F16" = 36
    "F17" = 4
    "F19" = 387
    "F20" = 10854
    "F21" = 10775
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
